# 更新工作区配置，替换旧文件 76. Minimum Window Substring 为新文件
# 3201. Find the Maximum Length of Valid Subsequence I，并新增
# 209. Minimum Size Subarray Sum 两条刷题记录。

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 84: 3201. Find the Maximum Length of Valid Subsequence I ----
$ws.Cells.Item(84, 1).Value = 3201
$ws.Cells.Item(84, 2).Value = "Find the Maximum Length of Valid Subsequence I"
$ws.Cells.Item(84, 3).Value = "#array #dynamic-programming "
$ws.Cells.Item(84, 4).Value = "medium"
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 1
$ws.Cells.Item(84, 7).Value = 50
$ws.Cells.Item(84, 8).Value = 45854
$ws.Cells.Item(84, 9).Value = 45854
$ws.Cells.Item(84, 10).Value = "?"

# Reuse the existing date number format (style index 4) from the row above
# instead of letting Excel mint a brand-new custom numFmt.
$ws.Range("H83:I83").Copy()
$ws.Range("H84:I84").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(84).RowHeight = 51

# ---- Row 85: 209. Minimum Size Subarray Sum ----
$ws.Cells.Item(85, 1).Value = 209
$ws.Cells.Item(85, 2).Value = "Minimum Size Subarray Sum"
$ws.Cells.Item(85, 3).Value = "#two-pointers #sliding-window #核心 "
$ws.Cells.Item(85, 4).Value = "medium"
$ws.Cells.Item(85, 5).Value = 2
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 30
$ws.Cells.Item(85, 8).Value = 44366
$ws.Cells.Item(85, 9).Value = 45854
# (no value in column J for row 85)

$ws.Range("H83:I83").Copy()
$ws.Range("H85:I85").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(85).RowHeight = 51

$excel.CutCopyMode = $false

# ---- Update the view: scroll/selection now centres on the new rows ----
$ws.Range("G85").Select() | Out-Null
